$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ستون اول"
$ws.Range("B1").Value = "ستون دوم"
$ws.Range("C1").Value = "ستون سوم"

$ws.Range("C2").Select()
